$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 410-437 with revised field values ---
# Row 410
$ws.Range("D410").Value = 44516
$ws.Range("I410").Value = 'Primera'
$ws.Range("J410").Value = 1000
$ws.Range("K410").Value = 22000
$ws.Range("M410").Value = 22500
$ws.Range("P410").Value = 1250

# Row 411
$ws.Range("D411").Value = 44516
$ws.Range("I411").Value = 'Segunda'
$ws.Range("J411").Value = 300
$ws.Range("K411").Value = 15000
$ws.Range("L411").Value = 15000
$ws.Range("M411").Value = 15000
$ws.Range("P411").Value = 833

# Row 412
$ws.Range("I412").Value = 'Extra'
$ws.Range("J412").Value = 600
$ws.Range("K412").Value = 23000
$ws.Range("L412").Value = 23000
$ws.Range("M412").Value = 23000
$ws.Range("P412").Value = 1278

# Row 413
$ws.Range("D413").Value = 44295
$ws.Range("I413").Value = 'Primera'
$ws.Range("J413").Value = 600
$ws.Range("K413").Value = 19000
$ws.Range("L413").Value = 19500
$ws.Range("M413").Value = 19250
$ws.Range("P413").Value = 1069

# Row 414
$ws.Range("D414").Value = 44295
$ws.Range("H414").Value = 'Larga vida'
$ws.Range("I414").Value = 'Segunda'
$ws.Range("J414").Value = 300
$ws.Range("K414").Value = 15000
$ws.Range("L414").Value = 15000
$ws.Range("M414").Value = 15000
$ws.Range("O414").Value = 'Limache'
$ws.Range("P414").Value = 833

# Row 415
$ws.Range("D415").Value = 44217
$ws.Range("I415").Value = 'Extra'
$ws.Range("J415").Value = 350
$ws.Range("K415").Value = 18000
$ws.Range("L415").Value = 18000
$ws.Range("M415").Value = 18000
$ws.Range("P415").Value = 1000

# Row 416
$ws.Range("D416").Value = 44217
$ws.Range("H416").Value = 'Semiduro'
$ws.Range("J416").Value = 100
$ws.Range("K416").Value = 14000
$ws.Range("L416").Value = 14000
$ws.Range("M416").Value = 14000
$ws.Range("N416").Value = '$/bandeja 18 kilos'
$ws.Range("O416").Value = 'Región de O''Higgins'
$ws.Range("P416").Value = 778
$ws.Range("Q416").Value = 18

# Row 417
$ws.Range("D417").Value = 44509
$ws.Range("I417").Value = 'Primera'
$ws.Range("J417").Value = 800
$ws.Range("K417").Value = 24000
$ws.Range("L417").Value = 25000
$ws.Range("M417").Value = 24500
$ws.Range("O417").Value = 'Limache'
$ws.Range("P417").Value = 1361

# Row 418
$ws.Range("D418").Value = 44509
$ws.Range("I418").Value = 'Primera'
$ws.Range("J418").Value = 500
$ws.Range("K418").Value = 24000
$ws.Range("L418").Value = 24000
$ws.Range("M418").Value = 24000
$ws.Range("P418").Value = 1200

# Row 419
$ws.Range("I419").Value = 'Extra'
$ws.Range("J419").Value = 250
$ws.Range("K419").Value = 18000
$ws.Range("L419").Value = 18000
$ws.Range("M419").Value = 18000
$ws.Range("P419").Value = 1000

# Row 420
$ws.Range("I420").Value = 'Extra'
$ws.Range("K420").Value = 20000
$ws.Range("L420").Value = 20000
$ws.Range("M420").Value = 20000
$ws.Range("P420").Value = 1000

# Row 421
$ws.Range("I421").Value = 'Primera'
$ws.Range("J421").Value = 500
$ws.Range("K421").Value = 15000
$ws.Range("L421").Value = 16000
$ws.Range("M421").Value = 15500
$ws.Range("P421").Value = 861

# Row 422
$ws.Range("D422").Value = 44421
$ws.Range("I422").Value = 'Primera'
$ws.Range("J422").Value = 300
$ws.Range("K422").Value = 19000
$ws.Range("L422").Value = 19000
$ws.Range("M422").Value = 19000
$ws.Range("N422").Value = '$/bandeja 20 kilos'
$ws.Range("P422").Value = 950
$ws.Range("Q422").Value = 20

# Row 423
$ws.Range("D423").Value = 44421
$ws.Range("I423").Value = 'Segunda'
$ws.Range("J423").Value = 250
$ws.Range("K423").Value = 14000
$ws.Range("L423").Value = 14000
$ws.Range("M423").Value = 14000
$ws.Range("N423").Value = '$/bandeja 18 kilos'
$ws.Range("P423").Value = 778
$ws.Range("Q423").Value = 18

# Row 424
$ws.Range("I424").Value = 'Extra'
$ws.Range("J424").Value = 500
$ws.Range("K424").Value = 17000
$ws.Range("L424").Value = 17500
$ws.Range("M424").Value = 17250
$ws.Range("P424").Value = 958

# Row 425
$ws.Range("I425").Value = 'Extra'
$ws.Range("J425").Value = 450
$ws.Range("K425").Value = 19500
$ws.Range("L425").Value = 20000
$ws.Range("M425").Value = 19778
$ws.Range("N425").Value = '$/bandeja 20 kilos'
$ws.Range("P425").Value = 989
$ws.Range("Q425").Value = 20

# Row 426
$ws.Range("D426").Value = 44383
$ws.Range("K426").Value = 15000
$ws.Range("L426").Value = 15000
$ws.Range("M426").Value = 15000
$ws.Range("O426").Value = 'Región de Arica y Parinacota'
$ws.Range("P426").Value = 833

# Row 427
$ws.Range("D427").Value = 44383
$ws.Range("I427").Value = 'Segunda'
$ws.Range("J427").Value = 250
$ws.Range("K427").Value = 14000
$ws.Range("L427").Value = 14000
$ws.Range("M427").Value = 14000
$ws.Range("O427").Value = 'Región de Arica y Parinacota'
$ws.Range("P427").Value = 778

# Row 428
$ws.Range("D428").Value = 44244
$ws.Range("J428").Value = 250
$ws.Range("K428").Value = 17500
$ws.Range("M428").Value = 17700
$ws.Range("P428").Value = 983

# Row 429
$ws.Range("D429").Value = 44273
$ws.Range("I429").Value = 'Extra'
$ws.Range("J429").Value = 150
$ws.Range("K429").Value = 19000
$ws.Range("L429").Value = 20000
$ws.Range("M429").Value = 19333
$ws.Range("O429").Value = 'Limache'
$ws.Range("P429").Value = 1074

# Row 430
$ws.Range("D430").Value = 44273
$ws.Range("J430").Value = 100
$ws.Range("K430").Value = 18000
$ws.Range("L430").Value = 18000
$ws.Range("M430").Value = 18000
$ws.Range("P430").Value = 1000

# Row 431
$ws.Range("D431").Value = 44433
$ws.Range("I431").Value = 'Primera'
$ws.Range("J431").Value = 120
$ws.Range("K431").Value = 15000
$ws.Range("L431").Value = 15000
$ws.Range("M431").Value = 15000
$ws.Range("O431").Value = 'Región de Arica y Parinacota'
$ws.Range("P431").Value = 833

# Row 432
$ws.Range("D432").Value = 44302
$ws.Range("J432").Value = 750
$ws.Range("K432").Value = 17000
$ws.Range("L432").Value = 17000
$ws.Range("M432").Value = 17000
$ws.Range("P432").Value = 944

# Row 433
$ws.Range("D433").Value = 44302
$ws.Range("I433").Value = 'Segunda'
$ws.Range("O433").Value = 'Limache'

# Row 434
$ws.Range("D434").Value = 44179
$ws.Range("I434").Value = 'Primera'
$ws.Range("J434").Value = 250
$ws.Range("K434").Value = 15000
$ws.Range("L434").Value = 15000
$ws.Range("M434").Value = 15000
$ws.Range("P434").Value = 833

# Row 435
$ws.Range("D435").Value = 44179
$ws.Range("I435").Value = 'Primera'
$ws.Range("J435").Value = 500
$ws.Range("K435").Value = 12000
$ws.Range("L435").Value = 12000
$ws.Range("M435").Value = 12000
$ws.Range("N435").Value = '$/bandeja 18 kilos'
$ws.Range("O435").Value = 'Región de Arica y Parinacota'
$ws.Range("P435").Value = 667
$ws.Range("Q435").Value = 18

# Row 436
$ws.Range("I436").Value = 'Extra'
$ws.Range("J436").Value = 300
$ws.Range("K436").Value = 24000
$ws.Range("L436").Value = 24000
$ws.Range("M436").Value = 24000
$ws.Range("O436").Value = 'Limache'
$ws.Range("P436").Value = 1333

# Row 437
$ws.Range("I437").Value = 'Extra'
$ws.Range("J437").Value = 300
$ws.Range("K437").Value = 25000
$ws.Range("L437").Value = 25000
$ws.Range("M437").Value = 25000
$ws.Range("N437").Value = '$/bandeja 20 kilos'
$ws.Range("O437").Value = 'Limache'
$ws.Range("P437").Value = 1250
$ws.Range("Q437").Value = 20

# --- Append new rows 438-439 ---
# Row 438
$ws.Range("D438").NumberFormat = $ws.Range("D410").NumberFormat()
$ws.Range("A438").Value = 4
$ws.Range("B438").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C438").Value = 'Los Lagos'
$ws.Range("D438").Value = 44491
$ws.Range("E438").Value = 10
$ws.Range("F438").Value = 100112020
$ws.Range("G438").Value = 'Tomate'
$ws.Range("H438").Value = 'Larga vida'
$ws.Range("I438").Value = 'Primera'
$ws.Range("J438").Value = 600
$ws.Range("K438").Value = 22000
$ws.Range("L438").Value = 22000
$ws.Range("M438").Value = 22000
$ws.Range("N438").Value = '$/bandeja 18 kilos'
$ws.Range("O438").Value = 'Región de Arica y Parinacota'
$ws.Range("P438").Value = 1222
$ws.Range("Q438").Value = 18
$ws.Range("R438").Value = 'Hortaliza'

# Row 439
$ws.Range("D439").NumberFormat = $ws.Range("D410").NumberFormat()
$ws.Range("A439").Value = 4
$ws.Range("B439").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C439").Value = 'Los Lagos'
$ws.Range("D439").Value = 44491
$ws.Range("E439").Value = 10
$ws.Range("F439").Value = 100112020
$ws.Range("G439").Value = 'Tomate'
$ws.Range("H439").Value = 'Larga vida'
$ws.Range("I439").Value = 'Primera'
$ws.Range("J439").Value = 500
$ws.Range("K439").Value = 9000
$ws.Range("L439").Value = 9000
$ws.Range("M439").Value = 9000
$ws.Range("N439").Value = '$/caja 10 kilos'
$ws.Range("O439").Value = 'Región de Arica y Parinacota'
$ws.Range("P439").Value = 900
$ws.Range("Q439").Value = 10
$ws.Range("R439").Value = 'Hortaliza'

